$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "78% Introverted"
$ws.Range("D2").Value = "77% Intuitive"
$ws.Range("E2").Value = "58% Feeling"
$ws.Range("F2").Value = "76% Judging"
$ws.Range("G2").Value = "63% Turbulent"
